# edit.ps1 - apply the "second commit" changes to Project Scope.docx
#
# 1) Remove the stray _GoBack bookmark that currently sits in the very
#    first paragraph (right after the manual line break).
# 2) Add a new centered / bold paragraph reading "Unique Poudel" right
#    after the "Squad 2-C2" paragraph, and move the _GoBack bookmark
#    (collapsed) to sit right after that new run.

$d = $word.ActiveDocument

# --- 1) drop the old _GoBack bookmark -------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2) locate the "Squad 2-C2" paragraph ----------------------------------
$squadParagraph = $null
$squadIndex = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "*Squad 2-C2*") {
        $squadParagraph = $p
        $squadIndex = $i
        break
    }
}

# Insert a brand-new paragraph straight after it; InsertParagraphAfter
# clones the surrounding centered/bold paragraph formatting automatically.
$squadParagraph.Range.InsertParagraphAfter()

$newParagraph = $d.Paragraphs.Item($squadIndex + 1)
$newParagraph.Range.InsertAfter("Unique Poudel")

$newRange = $newParagraph.Range
$bookmarkPos = $newRange.End - 1

# Collapsed (zero-length) ranges aren't positioned reliably by
# Bookmarks.Add in this host, so insert a throwaway placeholder
# character, anchor the bookmark around it (a real, non-empty range),
# then delete the placeholder. The bookmark collapses in place, right
# after "Unique Poudel" and before the paragraph mark - exactly like a
# normal Word _GoBack bookmark.
$placeholder = $d.Range($bookmarkPos, $bookmarkPos)
$placeholder.InsertAfter("X")

$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos + 1)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$cleanupRange = $d.Range($bookmarkPos, $bookmarkPos + 1)
$cleanupRange.Text = ""
